# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.517.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.003.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.599"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.371"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0747"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.69%  "
$ws.Range("E12").Value = "  -4.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.299.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.758"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.009.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.547.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.37%  "
$ws.Range("E21").Value = "  -4.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "221.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.92%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  -8.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("E28").Value = "  -2.99%  "
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "18.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.22%  "
$ws.Range("E32").Value = "  -2.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.91%  "
$ws.Range("E34").Value = "  -6.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.26%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.41%  "
$ws.Range("E41").Value = "  -1.94%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.452.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.50%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0925"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("E44").Value = "  -4.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("E46").Value = "  -8.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.65%  "
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +23.97%  "
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("E51").Value = "  -2.89%  "
